# Applies the row-content permutation edits described in the commit diff.
# Generated from a cell-by-cell diff of the worksheet "Artfynd".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Cells.Item(32, 1).Value = 111748730
$ws.Cells.Item(32, 2).Value = 90674
$ws.Cells.Item(32, 4).Value = 'VU'
$ws.Cells.Item(32, 5).Value = 2058
$ws.Cells.Item(32, 6).Value = 'Koppartaggsvamp'
$ws.Cells.Item(32, 7).Value = 'Hydnellum lundellii'
$ws.Cells.Item(32, 8).Value = '(Maas Geest. & Nannf.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Cells.Item(32, 17).Value = 664867.8416141088
$ws.Cells.Item(32, 18).Value = 6699418.15535864

# Row 33
$ws.Cells.Item(33, 1).Value = 111748691
$ws.Cells.Item(33, 2).Value = 89183
$ws.Cells.Item(33, 5).Value = 3215
$ws.Cells.Item(33, 6).Value = 'Rödgul trumpetsvamp'
$ws.Cells.Item(33, 7).Value = 'Craterellus lutescens'
$ws.Cells.Item(33, 8).Value = '(Fr.) Fr.'
$ws.Cells.Item(33, 17).Value = 664620.6444729717
$ws.Cells.Item(33, 18).Value = 6699349.02459933

# Row 34
$ws.Cells.Item(34, 1).Value = 111748644
$ws.Cells.Item(34, 2).Value = 99413
$ws.Cells.Item(34, 4).Value = 'LC'
$ws.Cells.Item(34, 5).Value = 221235
$ws.Cells.Item(34, 6).Value = 'Vårärt'
$ws.Cells.Item(34, 7).Value = 'Lathyrus vernus'
$ws.Cells.Item(34, 8).Value = '(L.) Bernh.'
$ws.Cells.Item(34, 17).Value = 664676.080560883
$ws.Cells.Item(34, 18).Value = 6699066.958814406
$ws.Cells.Item(34, 25).Value = "'2023-07-01"
$ws.Cells.Item(34, 27).Value = "'2023-07-01"

# Row 35
$ws.Cells.Item(35, 1).Value = 111748736
$ws.Cells.Item(35, 2).Value = 88966
$ws.Cells.Item(35, 4).Value = 'NT'
$ws.Cells.Item(35, 5).Value = 5754
$ws.Cells.Item(35, 6).Value = 'Gultoppig fingersvamp'
$ws.Cells.Item(35, 7).Value = 'Ramaria testaceoflava'
$ws.Cells.Item(35, 8).Value = '(Bres.) Corner'
$ws.Cells.Item(35, 17).Value = 664803.3322592583
$ws.Cells.Item(35, 18).Value = 6699212.309487404

# Row 36
$ws.Cells.Item(36, 1).Value = 111748704
$ws.Cells.Item(36, 2).Value = 85089
$ws.Cells.Item(36, 5).Value = 3762
$ws.Cells.Item(36, 6).Value = 'Olivspindling'
$ws.Cells.Item(36, 7).Value = 'Cortinarius venetus'
$ws.Cells.Item(36, 8).Value = '(Fr.:Fr.) Fr.'
$ws.Cells.Item(36, 17).Value = 664771.2603847764
$ws.Cells.Item(36, 18).Value = 6698957.93638737
$ws.Cells.Item(36, 25).Value = "'2023-08-26"
$ws.Cells.Item(36, 27).Value = "'2023-08-26"

# Row 40
$ws.Cells.Item(40, 1).Value = 111748683
$ws.Cells.Item(40, 2).Value = 96266
$ws.Cells.Item(40, 4).Value = 'LC'
$ws.Cells.Item(40, 5).Value = 223591
$ws.Cells.Item(40, 6).Value = 'Skogsnycklar'
$ws.Cells.Item(40, 7).Value = 'Dactylorhiza maculata subsp. fuchsii'
$ws.Cells.Item(40, 8).Value = '(Druce) Hyl.'
$ws.Cells.Item(40, 17).Value = 664638.5631493796
$ws.Cells.Item(40, 18).Value = 6699489.398619195

# Row 41
$ws.Cells.Item(41, 1).Value = 111748737
$ws.Cells.Item(41, 2).Value = 88918
$ws.Cells.Item(41, 4).Value = 'VU'
$ws.Cells.Item(41, 5).Value = 5745
$ws.Cells.Item(41, 6).Value = 'Gyllenfingersvamp'
$ws.Cells.Item(41, 7).Value = 'Ramaria brunneicontusa'
$ws.Cells.Item(41, 8).Value = 'R.H.Petersen'
$ws.Cells.Item(41, 17).Value = 664808.6128849701
$ws.Cells.Item(41, 18).Value = 6699172.460106185
$ws.Cells.Item(41, 25).Value = "'2023-08-26"
$ws.Cells.Item(41, 27).Value = "'2023-08-26"

# Row 42
$ws.Cells.Item(42, 1).Value = 111748649
$ws.Cells.Item(42, 2).Value = 94134
$ws.Cells.Item(42, 4).Value = 'NT'
$ws.Cells.Item(42, 5).Value = 53
$ws.Cells.Item(42, 6).Value = 'Vedtrappmossa'
$ws.Cells.Item(42, 7).Value = 'Crossocalyx hellerianus'
$ws.Cells.Item(42, 8).Value = '(Nees ex Lindenb.) Meyl.'
$ws.Cells.Item(42, 17).Value = 664603.5764787464
$ws.Cells.Item(42, 18).Value = 6699342.806076139
$ws.Cells.Item(42, 25).Value = "'2023-07-01"
$ws.Cells.Item(42, 27).Value = "'2023-07-01"

# Row 47
$ws.Cells.Item(47, 1).Value = 111748729
$ws.Cells.Item(47, 2).Value = 90674
$ws.Cells.Item(47, 4).Value = 'VU'
$ws.Cells.Item(47, 5).Value = 2058
$ws.Cells.Item(47, 6).Value = 'Koppartaggsvamp'
$ws.Cells.Item(47, 7).Value = 'Hydnellum lundellii'
$ws.Cells.Item(47, 8).Value = '(Maas Geest. & Nannf.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Cells.Item(47, 17).Value = 664589.6062576605
$ws.Cells.Item(47, 18).Value = 6699017.004385524

# Row 48
$ws.Cells.Item(48, 1).Value = 111748707
$ws.Cells.Item(48, 2).Value = 96369
$ws.Cells.Item(48, 5).Value = 219862
$ws.Cells.Item(48, 6).Value = 'Nästrot'
$ws.Cells.Item(48, 7).Value = 'Neottia nidus-avis'
$ws.Cells.Item(48, 8).Value = '(L.) Rich.'
$ws.Cells.Item(48, 17).Value = 664850.52293942
$ws.Cells.Item(48, 18).Value = 6699362.928853855

# Row 49
$ws.Cells.Item(49, 1).Value = 111748689
$ws.Cells.Item(49, 2).Value = 90678
$ws.Cells.Item(49, 4).Value = 'LC'
$ws.Cells.Item(49, 5).Value = 4366
$ws.Cells.Item(49, 6).Value = 'Skarp dropptaggsvamp'
$ws.Cells.Item(49, 7).Value = 'Hydnellum peckii'
$ws.Cells.Item(49, 8).Value = 'Banker'
$ws.Cells.Item(49, 17).Value = 664821.6049801367
$ws.Cells.Item(49, 18).Value = 6699355.67585695

# Row 60
$ws.Cells.Item(60, 1).Value = 111748739
$ws.Cells.Item(60, 2).Value = 103288
$ws.Cells.Item(60, 4).Value = 'LC'
$ws.Cells.Item(60, 5).Value = 221144
$ws.Cells.Item(60, 6).Value = 'Grönpyrola'
$ws.Cells.Item(60, 7).Value = 'Pyrola chlorantha'
$ws.Cells.Item(60, 8).Value = 'Sw.'
$ws.Cells.Item(60, 17).Value = 664649.5501774848
$ws.Cells.Item(60, 18).Value = 6699105.349934292

# Row 61
$ws.Cells.Item(61, 1).Value = 111748735
$ws.Cells.Item(61, 2).Value = 88966
$ws.Cells.Item(61, 5).Value = 5754
$ws.Cells.Item(61, 6).Value = 'Gultoppig fingersvamp'
$ws.Cells.Item(61, 7).Value = 'Ramaria testaceoflava'
$ws.Cells.Item(61, 8).Value = '(Bres.) Corner'
$ws.Cells.Item(61, 17).Value = 664801.9547313601
$ws.Cells.Item(61, 18).Value = 6699384.477357466
$ws.Cells.Item(61, 25).Value = "'2023-08-26"
$ws.Cells.Item(61, 27).Value = "'2023-08-26"

# Row 62
$ws.Cells.Item(62, 1).Value = 111748648
$ws.Cells.Item(62, 2).Value = 94134
$ws.Cells.Item(62, 5).Value = 53
$ws.Cells.Item(62, 6).Value = 'Vedtrappmossa'
$ws.Cells.Item(62, 7).Value = 'Crossocalyx hellerianus'
$ws.Cells.Item(62, 8).Value = '(Nees ex Lindenb.) Meyl.'
$ws.Cells.Item(62, 13).Value = ''
$ws.Cells.Item(62, 17).Value = 664655.3316318352
$ws.Cells.Item(62, 18).Value = 6699447.599252445

# Row 63
$ws.Cells.Item(63, 1).Value = 111748700
$ws.Cells.Item(63, 2).Value = 85089
$ws.Cells.Item(63, 5).Value = 3762
$ws.Cells.Item(63, 6).Value = 'Olivspindling'
$ws.Cells.Item(63, 7).Value = 'Cortinarius venetus'
$ws.Cells.Item(63, 8).Value = '(Fr.:Fr.) Fr.'
$ws.Cells.Item(63, 17).Value = 664853.490245176
$ws.Cells.Item(63, 18).Value = 6699352.175737353

# Row 64
$ws.Cells.Item(64, 1).Value = 111748705
$ws.Cells.Item(64, 2).Value = 96369
$ws.Cells.Item(64, 5).Value = 219862
$ws.Cells.Item(64, 6).Value = 'Nästrot'
$ws.Cells.Item(64, 7).Value = 'Neottia nidus-avis'
$ws.Cells.Item(64, 8).Value = '(L.) Rich.'
$ws.Cells.Item(64, 17).Value = 664770.1427733348
$ws.Cells.Item(64, 18).Value = 6699440.935377424

# Row 65
$ws.Cells.Item(65, 1).Value = 111748677
$ws.Cells.Item(65, 2).Value = 56414
$ws.Cells.Item(65, 4).Value = 'NT'
$ws.Cells.Item(65, 5).Value = 100049
$ws.Cells.Item(65, 6).Value = 'Spillkråka'
$ws.Cells.Item(65, 7).Value = 'Dryocopus martius'
$ws.Cells.Item(65, 8).Value = '(Linnaeus, 1758)'
$ws.Cells.Item(65, 13).Value = 'äldre spår'
$ws.Cells.Item(65, 17).Value = 664615.8542000444
$ws.Cells.Item(65, 18).Value = 6699498.265643556
$ws.Cells.Item(65, 25).Value = "'2023-07-01"
$ws.Cells.Item(65, 27).Value = "'2023-07-01"
